$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of Excel row number -> [new DAMSLTag, new DialogAct]
$changes = @{
    2   = @('sv', 'Statement-opinion')
    15  = @('b',  'Acknowledge (Backchannel)')
    23  = @('%',  'Uninterpretable')
    27  = @('aa', 'Agree/Accept')
    31  = @('b',  'Acknowledge (Backchannel)')
    59  = @('b',  'Acknowledge (Backchannel)')
    60  = @('sv', 'Statement-opinion')
    71  = @('aa', 'Agree/Accept')
    73  = @('aa', 'Agree/Accept')
    84  = @('sd', 'Statement-non-opinion')
    93  = @('sd', 'Statement-non-opinion')
    108 = @('sd', 'Statement-non-opinion')
    121 = @('aa', 'Agree/Accept')
    139 = @('sv', 'Statement-opinion')
    141 = @('sd', 'Statement-non-opinion')
    149 = @('sd', 'Statement-non-opinion')
    150 = @('sd', 'Statement-non-opinion')
    154 = @('b',  'Acknowledge (Backchannel)')
    168 = @('%',  'Uninterpretable')
    195 = @('sd', 'Statement-non-opinion')
    198 = @('sv', 'Statement-opinion')
}

foreach ($row in $changes.Keys) {
    $values = $changes[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}
